$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the used range to determine the last row with data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds the "Förändrad" (changed) date. All rows from 2 to last row
# had the serial date value 45203 which should become 45205.
$ws.Range("C2:C$lastRow").Value = 45205
